# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 00:52"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 816385
$ws.Range("C4").Value = 23626
$ws.Range("D4").Value = 82693
$ws.Range("E4").Value = 688518
$ws.Range("G4").Value = 2660
$ws.Range("H4").Value = 45174

# Row 16: Canada
$ws.Range("B16").Value = 38422
$ws.Range("C16").Value = 1593
$ws.Range("E16").Value = 23446
$ws.Range("G16").Value = 143
$ws.Range("H16").Value = 1833

# Row 40: Chequia
$ws.Range("B40").Value = 7033
$ws.Range("C40").Value = 133
$ws.Range("E40").Value = 5079

# Row 157: Bahamas
$ws.Range("B157").Value = 65
$ws.Range("C157").Value = 5
$ws.Range("D157").Value = 12
